$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row above the header row (old row 5) to hold the new
#    "Betrag pro Kind" input line. This pushes the header/data/totals rows
#    down by one (old 5->6, 6->7, 7->8) and Excel auto-adjusts formulas.
# ---------------------------------------------------------------------------
$ws.Rows("4:4").Insert()

# New label cell A4, styled like the neighbouring "Erstellt am" label (A3).
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Betrag pro Kind"

# New input cell B4, styled like the neighbouring date-input cell (B3),
# but formatted as a number with 2 decimals instead of a date.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "{betragProKind}"
$ws.Range("B4").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 2. Update the defined names so they keep pointing at the (now shifted)
#    placeholder data row (old row 6 -> new row 7).
# ---------------------------------------------------------------------------
$wb.Names.Item("kinderBereitsVerrechnet").RefersTo = "=Data!`$D`$7"
$wb.Names.Item("kinderTotal").RefersTo = "=Data!`$C`$7"

# ---------------------------------------------------------------------------
# 3. Add the new "Betrag zu verrechnen" summary line in row 10 (row 9 left
#    blank for spacing, matching the sheet's existing visual rhythm).
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Betrag zu verrechnen"
$ws.Range("A10").Borders.Item(7).LineStyle = 1
$ws.Range("A10").Borders.Item(8).LineStyle = 1
$ws.Range("A10").Borders.Item(9).LineStyle = 1

$r = $ws.Range("B10:D10")
$r.Borders.Item(8).LineStyle = 1
$r.Borders.Item(9).LineStyle = 1

# Result cell: copy an existing fully-boxed style, then set formula + format.
$ws.Range("A7").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").NumberFormat = "0.00"
$ws.Range("E10").Formula = "=B4*E8"

$excel.CutCopyMode = 0

Write-Host "edit applied"
